# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For each data row, take the comma-separated list of recorders in column G.
# Any entries that are exactly "System" are moved to the front (keeping their
# relative order/count), while the remaining entries keep their original
# relative order. If there is no "System" entry, the list is simply reversed.
# Rows with a single value are unaffected (no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "

    $systemParts = @($parts | Where-Object { $_.Equals("System") })
    $restParts = @($parts | Where-Object { -not ($_.Equals("System")) })

    if ($systemParts.Count -gt 0) {
        $newParts = @($systemParts) + @($restParts)
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newText = [string]::Join(", ", $newParts)

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
